$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-17 Friday" "2024-05-18 Saturday"

Replace-Text "99×26=2574" "79×22=1738"
Replace-Text "39×89=3471" "36×45=1620"
Replace-Text "71×82=5822" "25×81=2025"
Replace-Text "71×36=2556" "78×11=858"
Replace-Text "27×12=324" "73×19=1387"

Replace-Text "90×22=1980" "99×76=7524"
Replace-Text "15×82=1230" "89×47=4183"
Replace-Text "19×43=817" "85×85=7225"
Replace-Text "84×62=5208" "96×95=9120"
Replace-Text "34×20=680" "94×71=6674"

Replace-Text "22×63=1386" "76×75=5700"
Replace-Text "41×20=820" "49×25=1225"
Replace-Text "66×37=2442" "92×22=2024"
Replace-Text "35×22=770" "50×25=1250"
Replace-Text "87×60=5220" "77×29=2233"

Replace-Text "39×22=858" "71×44=3124"
Replace-Text "44×90=3960" "97×31=3007"
Replace-Text "53×11=583" "50×37=1850"
Replace-Text "26×64=1664" "53×83=4399"
Replace-Text "11×57=627" "16×31=496"

Replace-Text "21×85=1785" "16×29=464"
Replace-Text "86×27=2322" "64×78=4992"
Replace-Text "50×81=4050" "21×66=1386"
Replace-Text "42×35=1470" "16×55=880"
Replace-Text "12×36=432" "78×44=3432"
